$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("M1").Value = "fringilla"
$ws.Range("M5").Value = "C:\Users\Mark\Dropbox\LivemRNAData"
$ws.Range("M9").Value = "C:\Users\Mark\GregorLab\mRNADynamics"
$ws.Range("M4").Value = "C:\Users\Mark\GregorLab\Data\FISHAnalysisData"
$ws.Range("M3").Value = "C:\Users\Mark\GregorLab\Data\RawData"

$ws.Range("M3").Select() | Out-Null
$excel.ActiveWindow.SplitColumn = 10
